# Auto-generated PowerShell Excel COM-interop script
# Applies daily market-data refresh to the cryptos worksheet:
#  - Updates Price (D) and Volume/1h change (E) values for most rows
#  - Rows 42/43 swap coin identity (Aave <-> InjectiveProtocol) with new data
#
# Price/percentage values are stored as text in this sheet (even when they
# look numeric), so numeric-looking values are written with a temporary "@"
# (Text) number format to stop Excel from auto-converting them to real numbers,
# then the cell style is reset back to Normal to avoid leaving formatting behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.806.97"
$ws.Range("E2").Value = "  -1.00%  "
$ws.Range("D3").Value = "2.033.08"
$ws.Range("E3").Value = "  -1.53%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "227.16"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.24%  "
$ws.Range("E6").Value = "  -0.59%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "60.00"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.89%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.387"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.31%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0818"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.05%  "
$ws.Range("E11").Value = "  +0.08%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "14.65"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.09%  "
$ws.Range("D13").Value = "2.332.68"
$ws.Range("E13").Value = "  -1.58%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "21.07"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.34%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.760"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.71%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.22"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.21%  "
$ws.Range("D17").Value = "2.044.72"
$ws.Range("E17").Value = "  -1.26%  "
$ws.Range("D18").Value = "37.748.29"
$ws.Range("E18").Value = "  -0.92%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.07"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.31%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "69.83"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.11%  "
$ws.Range("E21").Value = "  -1.09%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "225.60"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.29%  "
$ws.Range("E23").Value = "  +0.13%  "
$ws.Range("E24").Value = "  -2.49%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.22"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.64%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.25"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.39%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "165.20"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.24%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.129"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -4.31%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "18.88"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.16%  "
$ws.Range("E30").Value = "  -7.26%  "
$ws.Range("E31").Value = "  +1.08%  "
$ws.Range("E32").Value = "  -2.92%  "
$ws.Range("E33").Value = "  +3.61%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0602"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.08%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.47"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.05%  "
$ws.Range("E36").Value = "  +6.24%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.25"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -5.63%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.24"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.66%  "
$ws.Range("E39").Value = "  +0.12%  "
$ws.Range("D40").Value = "1.541.41"
$ws.Range("E40").Value = "  +4.04%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0216"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.40%  "

# Rows 42 and 43 swap coins (Aave <-> InjectiveProtocol) with updated values
$ws.Range("B42").Value = "InjectiveProtocol"
$ws.Range("C42").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "16.89"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.12%  "
$ws.Range("B43").Value = "Aave"
$ws.Range("C43").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "96.78"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.73%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.83"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.65%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0921"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.59%  "
$ws.Range("E46").Value = "  -1.56%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.90"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -5.54%  "
$ws.Range("E48").Value = "  -1.98%  "
$ws.Range("E49").Value = "  -0.42%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.15"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.62%  "
$ws.Range("D51").Value = "2.223.16"
$ws.Range("E51").Value = "  -1.57%  "
